$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("http://sipp.pn-jakartapusat.go.id/", "https://sipp.pn-negara.go.id/", "https://sipp.pn-semarangkota.go.id/", "https://sipp.pn-surabayakota.go.id/")

for ($i = 0; $i -lt 2; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $row = 6 + ($i * 4) + $j
        $ws.Range("A$row").Value = $values[$j]
    }
}

$ws.Range("F9").Select()
